$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on cells whose new values would otherwise be
# auto-converted to numbers (e.g. "1.00", "228.36") so they remain
# stored as plain text, matching the source data (inline strings).
$textCells = @("D2", "D3", "D4", "D5", "D7", "D8", "D10", "D12", "D13", "D14", "D16", "D17", "D18", "D19", "D21", "D22", "D24", "D25", "D26", "D29", "D30", "D31", "D33", "D36", "D37", "D38", "D39", "D41", "D43", "D44", "D45", "D46", "D49", "D50")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '35.310.50'
$ws.Range('E2').Value = '  +2.09%  '
$ws.Range('D3').Value = '1.845.95'
$ws.Range('E3').Value = '  +2.04%  '
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  -0.10%  '
$ws.Range('D5').Value = '228.36'
$ws.Range('E5').Value = '  +1.19%  '
$ws.Range('E6').Value = '  +2.54%  '
$ws.Range('D7').Value = '0.999'
$ws.Range('E7').Value = '  -0.14%  '
$ws.Range('D8').Value = '43.12'
$ws.Range('E8').Value = '  +15.18%  '
$ws.Range('E9').Value = '  +4.87%  '
$ws.Range('D10').Value = '0.0695'
$ws.Range('E10').Value = '  +1.97%  '
$ws.Range('D12').Value = '2.111.15'
$ws.Range('E12').Value = '  +1.95%  '
$ws.Range('D13').Value = '11.64'
$ws.Range('E13').Value = '  +2.83%  '
$ws.Range('D14').Value = '1.837.73'
$ws.Range('E14').Value = '  +1.13%  '
$ws.Range('E15').Value = '  +7.17%  '
$ws.Range('D16').Value = '0.661'
$ws.Range('E16').Value = '  +4.56%  '
$ws.Range('D17').Value = '35.216.56'
$ws.Range('E17').Value = '  +1.95%  '
$ws.Range('D18').Value = '69.93'
$ws.Range('E18').Value = '  +1.85%  '
$ws.Range('D19').Value = '246.25'
$ws.Range('E19').Value = '  +1.14%  '
$ws.Range('E20').Value = '  +2.39%  '
$ws.Range('D21').Value = '12.13'
$ws.Range('E21').Value = '  +8.21%  '
$ws.Range('D22').Value = '4.73'
$ws.Range('E22').Value = '  +14.34%  '
$ws.Range('E23').Value = '  -0.05%  '
$ws.Range('D24').Value = '2.19'
$ws.Range('E24').Value = '  -1.51%  '
$ws.Range('D25').Value = '172.49'
$ws.Range('E25').Value = '  +0.14%  '
$ws.Range('D26').Value = '7.94'
$ws.Range('E26').Value = '  +1.04%  '
$ws.Range('E27').Value = '  +3.63%  '
$ws.Range('E28').Value = '  +1.62%  '
$ws.Range('D29').Value = '3.574.60'
$ws.Range('E29').Value = '  +47.12%  '
$ws.Range('D30').Value = '1.00'
$ws.Range('E30').Value = '  -0.01%  '
$ws.Range('D31').Value = '1.32'
$ws.Range('E31').Value = '  +7.99%  '
$ws.Range('E32').Value = '  +3.58%  '
$ws.Range('D33').Value = '4.06'
$ws.Range('E33').Value = '  +3.66%  '
$ws.Range('E34').Value = '  +4.02%  '
$ws.Range('E35').Value = '  +4.22%  '
$ws.Range('D36').Value = '0.675'
$ws.Range('E36').Value = '  +3.25%  '
$ws.Range('D37').Value = '90.53'
$ws.Range('E37').Value = '  +12.03%  '
$ws.Range('B38').Value = 'Maker'
$ws.Range('C38').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D38').Value = '1.344.01'
$ws.Range('E38').Value = '  -1.55%  '
$ws.Range('B39').Value = 'TrustWalletToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D39').Value = '1.08'
$ws.Range('E39').Value = '  +1.15%  '
$ws.Range('E40').Value = '  +9.16%  '
$ws.Range('D41').Value = '2.44'
$ws.Range('E41').Value = '  +2.98%  '
$ws.Range('E42').Value = '  +3.69%  '
$ws.Range('D43').Value = '14.85'
$ws.Range('E43').Value = '  +8.31%  '
$ws.Range('D44').Value = '1.26'
$ws.Range('E44').Value = '  +6.07%  '
$ws.Range('B45').Value = 'HuobiToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D45').Value = '2.45'
$ws.Range('E45').Value = '  +1.01%  '
$ws.Range('B46').Value = 'MXToken'
$ws.Range('C46').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D46').Value = '2.83'
$ws.Range('E46').Value = '  +1.38%  '
$ws.Range('E47').Value = '  +3.70%  '
$ws.Range('E48').Value = '  +4.34%  '
$ws.Range('D49').Value = '2.010.95'
$ws.Range('E49').Value = '  +2.03%  '
$ws.Range('D50').Value = '104.82'
$ws.Range('E50').Value = '  +1.97%  '
$ws.Range('E51').Value = '  -0.06%  '
